# Logic tree input file updated
# A new "Possible_Problem" node (with its percentage breakdown) was inserted
# as the new row 8, pushing the former rows 8-9 (the Yes/No branches) down
# to rows 9-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 - shifts old rows 8 & 9 down to 9 & 10.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the "Possible_Problem" node data
# (reusing the same text that already exists elsewhere in the workbook).
$ws.Range("A8").Value = "Problem:Does driving the vehicle alleviate the problem? (Please answer as: Yes, No)"
$ws.Range("B8").Value = "Possible_Problem"
$ws.Range("C8").Value = "Possible_Problem:30% Restricted Heater Core`n20%Thermostat`n15% HVAC Door Actuators`n10% Restricted Radiator Coolant Flow`n10% Low Coolant Level`n5% HVAC Control Unit`n5% Water Pump`n5% Cylinder Head Gasket"

# Match the tall row height used for this long wrapped note.
$ws.Rows.Item(8).RowHeight = 409.6

# Reflect the author's on-screen scroll position / active cell selection.
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("C8").Select()
